$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = 'I changed the API one more time for the flight data to open sky API But some of the imports werent working so I had to do some research and find where the JAR files are. 2 of the imports are still erroring, and I don’t know where to find the JAR file for it. '

$ws.Range("A47").Value = 45237
$ws.Range("A47").NumberFormat = "d-mmm"
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = $newText

$ws.Range("C47").Select()
